$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUuid = "7580109a-a8e0-4789-aaa6-bddca426b60e"
$newUuid = "c4f1636a-d92c-45e5-82c9-d146dc727e93"

# --- 1) Update existing rows 2-9: refresh uuid (col G) and recompute
#        downtime (col H) as ROUND(tgap * 60, 0) instead of tgap/1440 days.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 7).Value = $newUuid

    $tgap = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 8).Value = [Math]::Round($tgap * 60)
}

# --- 2) Append new rows 10-25 with the same column layout as the
#        existing data rows (A..H), reusing the number formats of row 2
#        for the date columns (D,E) and the downtime column (H).
$newRows = @(
    @(10, "Line:9 Stage:1", "01/09/2024", "pri cl LA", 45300.40887037037, 45300.40910185185, 0.33, 20),
    @(11, "Line:9 Stage:1", "01/09/2024", "Pri pH flows", 45300.49280315972, 45300.49893741898, 8.83, 530),
    @(12, "Line:9 Stage:1", "01/09/2024", "Pri pH flows", 45300.50241207176, 45300.50287503472, 0.67, 40),
    @(13, "Line:9 Stage:1", "01/09/2024", "Pri pH flows", 45300.50356979167, 45300.50368553241, 0.17, 10),
    @(14, "Line:9 Stage:1", "01/09/2024", "Pri pH flows", 45300.50438012731, 45300.50449586806, 0.17, 10),
    @(15, "Line:9 Stage:1", "01/09/2024", "Pri pH flows", 45300.50519054398, 45300.50715813658, 2.83, 170),
    @(16, "Line:3 Stage:1", "12/11/2023", "SEC Cl", 45271.36145825232, 45271.36158556713, 0.18, 11),
    @(17, "Line:3 Stage:1", "12/11/2023", "PRI pH", 45271.68537890046, 45271.68549464121, 0.17, 10),
    @(18, "Line:3 Stage:1", "12/11/2023", "PRI pH SEC pH", 45271.68549475694, 45271.68769383102, 3.17, 190),
    @(19, "Line:3 Stage:1", "12/11/2023", "PRI pH&rem SEC pH", 45271.6876965625, 45271.68898128472, 1.85, 111),
    @(20, "Line:3 Stage:1", "12/11/2023", "PRI pH&rem SEC pH&rem", 45271.68897129629, 45271.68908703703, 0.17, 10),
    @(21, "Line:3 Stage:1", "12/11/2023", "PRI pH SEC pH", 45271.6890871875, 45271.68920292824, 0.17, 10),
    @(22, "Line:3 Stage:1", "12/11/2023", "PRI pH SEC pH", 45271.89825825232, 45271.8996471412, 2, 120),
    @(23, "Line:3 Stage:1", "12/11/2023", "PRI Cl&pH", 45272.09804105324, 45272.09827253472, 0.33, 20),
    @(24, "Line:8 Stage:1", "01/09/2024", "pri cl LA", 45300.42056299769, 45300.42067873842, 0.17, 10),
    @(25, "Line:8 Stage:1", "01/09/2024", "pri pH HA", 45300.65476473379, 45300.65488047454, 0.17, 10)
)

$dateFormat = $ws.Range("D2").NumberFormat
$downtimeFormat = $ws.Range("H2").NumberFormat

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]

    # Column B (tdate) holds plain text like "01/09/2024" (not a real
    # date cell) in the source data. Force text so Excel doesn't
    # auto-convert the literal into a date serial, then drop back to the
    # workbook's default ("Normal") style so no stray number format is
    # left attached to the cell.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Value = $row[3]

    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 5).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 6).Value = $row[6]

    $ws.Cells.Item($r, 7).Value = $newUuid

    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 8).NumberFormat = $downtimeFormat
}
